$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The sheet is protected; unprotect it so the cells below can be updated, then
# re-protect it afterwards so the sheet stays protected like the original file.
$ws.Unprotect()

# Update the confidential disclaimer text (date change 2021-03-25 -> 2021-03-26)
$ws.Range("A16").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-26 for illustrative purposes only and are subject to change."
# Setting multi-line text can make Excel stamp an explicit row height; AutoFit
# brings row 16 back to the original "no explicit height" state.
$ws.Rows.Item(16).AutoFit()

# Update Weight (D) and Percent Change (E) values for rows 2-13
$ws.Range("D2").Value = 0.03051953815244842
$ws.Range("E2").Value = 0.02022058823529416

$ws.Range("D3").Value = 0.02430002239040803
$ws.Range("E3").Value = 0.00512070226773953

$ws.Range("D4").Value = 0.05277479687260463
$ws.Range("E4").Value = 0.01705653021442499

$ws.Range("D5").Value = 0.1380572891872484
$ws.Range("E5").Value = 0.01674145667932336

$ws.Range("D6").Value = 0.03128466983225348
$ws.Range("E6").Value = 0.02712886209495102

$ws.Range("D7").Value = 0.1198115020131196
$ws.Range("E7").Value = 0.0106908993719097

$ws.Range("D8").Value = 0.1016486702921646
$ws.Range("E8").Value = 0.01920438957476001

$ws.Range("D9").Value = 0.027924091473894
$ws.Range("E9").Value = 0.0260188809578632

$ws.Range("D10").Value = 0.1227185989656035
$ws.Range("E10").Value = 0.01623176494760625

$ws.Range("D11").Value = 0.2471782538709183
$ws.Range("E11").Value = 0.02354879348774119

$ws.Range("D12").Value = 0.1037825669493372
$ws.Range("E12").Value = -0.005918367346938846

$ws.Range("E13").Value = 0.01595972147139291

# Restore sheet protection (matches original protected state).
$ws.Protect()

$wb.Save()
